$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Update row 2: cpf becomes numeric 1 (value stays numeric)
$ws.Range("A2").Value = 1

# Update row 3: cpf becomes text "1" (not numeric) and profile becomes "professor"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1"
$ws.Range("A3").Style = "Normal"
$ws.Range("C3").Value = "professor"

# Remove rows 4 through 8 (old extra data)
$ws.Range("A4:C8").EntireRow.Delete()
